# Apply cryptos list update (GitHub Actions refresh) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.225.86"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "2.646.73"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Formula = "'598.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").Formula = "'156.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("E9").Value = "  +6.09%  "

$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("D13").Formula = "'28.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").Value = "3.129.40"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "68.243.43"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").Value = "2.659.81"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Formula = "'364.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("D20").Formula = "'7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").Formula = "'4.40"
$ws.Range("D21").Style = "Normal"

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("D24").Formula = "'75.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.68%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Formula = "'9.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.790.30"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Formula = "'0.0000105"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("D29").Formula = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Formula = "'556.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "

$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Formula = "'1.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("E34").Value = "  +1.98%  "

$ws.Range("D35").Formula = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("E36").Value = "  +3.36%  "

$ws.Range("D37").Formula = "'160.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").Formula = "'19.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").Formula = "'158.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").Formula = "'22.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("E51").Value = "  +0.86%  "
